$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.04412441691274
$ws.Range("D2").Value = 1.053205378258711
$ws.Range("E2").Value = 1.057909509656376
$ws.Range("F2").Value = 1.065005961470826
$ws.Range("I2").Value = 1.046080679249341
$ws.Range("J2").Value = 1.049191183350932
$ws.Range("K2").Value = 1.055951989628414
$ws.Range("L2").Value = 1.060643189174065
$ws.Range("M2").Value = 1.06772036649294
$ws.Range("N2").Value = 1.050681155175201

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.045007687974078
$ws.Range("D3").Value = 1.053935011560526
$ws.Range("E3").Value = 1.05884279773517
$ws.Range("F3").Value = 1.06592500943221
$ws.Range("I3").Value = 1.046334198122602
$ws.Range("J3").Value = 1.049721909684294
$ws.Range("K3").Value = 1.056494834988576
$ws.Range("L3").Value = 1.06139010337202
$ws.Range("M3").Value = 1.068454470075673
$ws.Range("N3").Value = 1.051212635200832

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.045579565676865
$ws.Range("D4").Value = 1.054407417981282
$ws.Range("E4").Value = 1.059447926233275
$ws.Range("F4").Value = 1.066520618901129
$ws.Range("I4").Value = 1.046497069572438
$ws.Range("J4").Value = 1.050065008084149
$ws.Range("K4").Value = 1.056845706469799
$ws.Range("L4").Value = 1.061874004772537
$ws.Range("M4").Value = 1.068929773769331
$ws.Range("N4").Value = 1.051556220839768

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.045820063759308
$ws.Range("D5").Value = 1.054606084350998
$ws.Range("E5").Value = 1.059702614881722
$ws.Range("F5").Value = 1.066771232770562
$ws.Range("I5").Value = 1.046565259420685
$ws.Range("J5").Value = 1.050209169714762
$ws.Range("K5").Value = 1.056993119361037
$ws.Range("L5").Value = 1.062077578780311
$ws.Range("M5").Value = 1.069129659440613
$ws.Range("N5").Value = 1.051700587196432

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.045860449141605
$ws.Range("D6").Value = 1.054639445177199
$ws.Range("E6").Value = 1.059745395341987
$ws.Range("F6").Value = 1.066813324777135
$ws.Range("I6").Value = 1.046576692280437
$ws.Range("J6").Value = 1.050233370545063
$ws.Range("K6").Value = 1.057017865119711
$ws.Range("L6").Value = 1.062111768033019
$ws.Range("M6").Value = 1.069163225073205
$ws.Range("N6").Value = 1.051724822394688

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.045582778908482
$ws.Range("D7").Value = 1.054410072310785
$ws.Range("E7").Value = 1.059451328247544
$ws.Range("F7").Value = 1.066523966754797
$ws.Range("I7").Value = 1.046497981834727
$ws.Range("J7").Value = 1.050066934682524
$ws.Range("K7").Value = 1.056847676576788
$ws.Range("L7").Value = 1.061876724382607
$ws.Range("M7").Value = 1.068932444385931
$ws.Range("N7").Value = 1.051558150174133

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.044422850887739
$ws.Range("D8").Value = 1.05345190136523
$ws.Range("E8").Value = 1.058224663402723
$ws.Range("F8").Value = 1.065316365780467
$ws.Range("I8").Value = 1.046166599445485
$ws.Range("J8").Value = 1.049370609817215
$ws.Range("K8").Value = 1.05613552595755
$ws.Range("L8").Value = 1.06089548773097
$ws.Range("M8").Value = 1.067968399475176
$ws.Range("N8").Value = 1.050860836447647

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.042381585726858
$ws.Range("D9").Value = 1.051765729465127
$ws.Range("E9").Value = 1.056072601119877
$ws.Range("F9").Value = 1.063195563886838
$ws.Range("I9").Value = 1.045573713582563
$ws.Range("J9").Value = 1.048141214304249
$ws.Range("K9").Value = 1.054877720803433
$ws.Range("L9").Value = 1.059171060420324
$ws.Range("M9").Value = 1.066271905806977
$ws.Range("N9").Value = 1.049629695052039

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.041022614673324
$ws.Range("D10").Value = 1.05064321441152
$ws.Range("E10").Value = 1.054644355862133
$ws.Range("F10").Value = 1.061786592136019
$ws.Range("I10").Value = 1.045172480461326
$ws.Range("J10").Value = 1.047320077062196
$ws.Range("K10").Value = 1.054037297582125
$ws.Range("L10").Value = 1.05802463671476
$ws.Range("M10").Value = 1.065142519392152
$ws.Range("N10").Value = 1.048807391700916

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.040434624406201
$ws.Range("D11").Value = 1.050157550863614
$ws.Range("E11").Value = 1.054027461763234
$ws.Range("F11").Value = 1.06117767236134
$ws.Range("I11").Value = 1.04499733300944
$ws.Range("J11").Value = 1.046964163017937
$ws.Range("K11").Value = 1.053672951476536
$ws.Range("L11").Value = 1.057528996345271
$ws.Range("M11").Value = 1.064653880750297
$ws.Range("N11").Value = 1.04845097221789

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.0402162879005
$ws.Range("D12").Value = 1.049977214294108
$ws.Range("E12").Value = 1.053798553541171
$ws.Range("F12").Value = 1.0609516704025
$ws.Range("I12").Value = 1.044932064119551
$ws.Range("J12").Value = 1.04683190818133
$ws.Range("K12").Value = 1.053537552500514
$ws.Range("L12").Value = 1.057345010204508
$ws.Range("M12").Value = 1.064472438966891
$ws.Range("N12").Value = 1.048318529564246

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.040263118669058
$ws.Range("D13").Value = 1.050015894339239
$ws.Range("E13").Value = 1.053847644531461
$ws.Range("F13").Value = 1.061000140512257
$ws.Range("I13").Value = 1.044946074070753
$ws.Range("J13").Value = 1.04686027965299
$ws.Range("K13").Value = 1.053566598946715
$ws.Range("L13").Value = 1.057384470552041
$ws.Range("M13").Value = 1.0645113560928
$ws.Range("N13").Value = 1.048346941326649

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.040416575222833
$ws.Range("D14").Value = 1.050142642943911
$ws.Range("E14").Value = 1.054008535361731
$ws.Range("F14").Value = 1.061158987329764
$ws.Range("I14").Value = 1.044991942169965
$ws.Range("J14").Value = 1.046953231857349
$ws.Range("K14").Value = 1.053661760666654
$ws.Range("L14").Value = 1.057513785604196
$ws.Range("M14").Value = 1.064638881476922
$ws.Range("N14").Value = 1.048440025533801

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.04051113408317
$ws.Range("D15").Value = 1.050220745002822
$ws.Range("E15").Value = 1.054107696527017
$ws.Range("F15").Value = 1.061256881723651
$ws.Range("I15").Value = 1.04502017503587
$ws.Range("J15").Value = 1.047010495849692
$ws.Range("K15").Value = 1.053720384414229
$ws.Range("L15").Value = 1.057593476385979
$ws.Range("M15").Value = 1.064717462112835
$ws.Range("N15").Value = 1.048497370847582

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.041061647298909
$ws.Range("D16").Value = 1.050675454704756
$ws.Range("E16").Value = 1.054685329783546
$ws.Range("F16").Value = 1.061827029007878
$ws.Range("I16").Value = 1.045184074750062
$ws.Range("J16").Value = 1.047343690485881
$ws.Range("K16").Value = 1.054061468912738
$ws.Range("L16").Value = 1.058057547037595
$ws.Range("M16").Value = 1.065174957123362
$ws.Range("N16").Value = 1.048831038658372

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.041407091849962
$ws.Range("D17").Value = 1.050960788310676
$ws.Range("E17").Value = 1.055048078714542
$ws.Range("F17").Value = 1.062184982809505
$ws.Range("I17").Value = 1.04528650743266
$ws.Range("J17").Value = 1.047552600020201
$ws.Range("K17").Value = 1.054275305897796
$ws.Range("L17").Value = 1.058348852990843
$ws.Range("M17").Value = 1.065462037826759
$ws.Range("N17").Value = 1.049040244868193

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.041608627457086
$ws.Range("D18").Value = 1.051127256383959
$ws.Range("E18").Value = 1.055259813047514
$ws.Range("F18").Value = 1.062393884405043
$ws.Range("I18").Value = 1.045346118512694
$ws.Range("J18").Value = 1.047674418981776
$ws.Range("K18").Value = 1.054399991126618
$ws.Range("L18").Value = 1.058518840903609
$ws.Range("M18").Value = 1.06562952498786
$ws.Range("N18").Value = 1.049162236826666

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.041677353258164
$ws.Range("D19").Value = 1.051184024063229
$ws.Range("E19").Value = 1.055332034213418
$ws.Range("F19").Value = 1.062465133604523
$ws.Range("I19").Value = 1.045366421234193
$ws.Range("J19").Value = 1.047715950228306
$ws.Range("K19").Value = 1.054442498365166
$ws.Range("L19").Value = 1.058576814915158
$ws.Range("M19").Value = 1.065686640160932
$ws.Range("N19").Value = 1.04920382705233

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.041370024395965
$ws.Range("D20").Value = 1.050930170793795
$ws.Range("E20").Value = 1.055009143762603
$ws.Range("F20").Value = 1.062146566034426
$ws.Range("I20").Value = 1.045275531455386
$ws.Range("J20").Value = 1.04753018955883
$ws.Range("K20").Value = 1.054252367576992
$ws.Range("L20").Value = 1.058317590947399
$ws.Range("M20").Value = 1.065431232860412
$ws.Range("N20").Value = 1.049017802581395

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.040371384171175
$ws.Range("D21").Value = 1.050105316957156
$ws.Range("E21").Value = 1.053961150558494
$ws.Range("F21").Value = 1.061112205976336
$ws.Range("I21").Value = 1.044978440986526
$ws.Range("J21").Value = 1.04692586119152
$ws.Range("K21").Value = 1.053633739686766
$ws.Range("L21").Value = 1.057475702310983
$ws.Range("M21").Value = 1.064601326752959
$ws.Range("N21").Value = 1.048412615998487

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.039743900726704
$ws.Range("D22").Value = 1.049587049352607
$ws.Range("E22").Value = 1.053303588582417
$ws.Range("F22").Value = 1.060462893122471
$ws.Range("I22").Value = 1.044790426057117
$ws.Range("J22").Value = 1.046545592630902
$ws.Range("K22").Value = 1.053244410506297
$ws.Range("L22").Value = 1.056947049687222
$ws.Range("M22").Value = 1.064079881742042
$ws.Range("N22").Value = 1.048031807412895

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.040076503140512
$ws.Range("D23").Value = 1.04986175907857
$ws.Range("E23").Value = 1.05365204587698
$ws.Range("F23").Value = 1.060807007972192
$ws.Range("I23").Value = 1.044890211980513
$ws.Range("J23").Value = 1.046747208607342
$ws.Range("K23").Value = 1.053450836288889
$ws.Range("L23").Value = 1.057227233924477
$ws.Range("M23").Value = 1.064356276007077
$ws.Range("N23").Value = 1.048233709707148

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.041386773448122
$ws.Range("D24").Value = 1.050944005415098
$ws.Range("E24").Value = 1.055026736331172
$ws.Range("F24").Value = 1.062163924571851
$ws.Range("I24").Value = 1.045280491447862
$ws.Range("J24").Value = 1.047540315988173
$ws.Range("K24").Value = 1.054262732547501
$ws.Range("L24").Value = 1.058331716690417
$ws.Range("M24").Value = 1.065445152181923
$ws.Range("N24").Value = 1.04902794339143

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.042908977038561
$ws.Range("D25").Value = 1.052201369849379
$ws.Range("E25").Value = 1.056627828303711
$ws.Range("F25").Value = 1.063742985772942
$ws.Range("I25").Value = 1.045728045020033
$ws.Range("J25").Value = 1.04845931848593
$ws.Range("K25").Value = 1.055203230960114
$ws.Range("L25").Value = 1.05961630879211
$ws.Range("M25").Value = 1.066710212189116
$ws.Range("N25").Value = 1.049948250978145
